$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 2, pushing the existing tasks down.
$ws.Range("A2:C3").EntireRow.Insert()

# New row 2: Engine / Complete Renderer refactoring / 10
$ws.Range("A2").Value = "Engine"
$ws.Range("B2").Value = "Complete Renderer refactoring"
$ws.Range("C2").Value = 10

# New row 3: Engine / Deprecate Vector3.  Replace it with position, direction and unit direction / 10
$ws.Range("A3").Value = "Engine"
$ws.Range("B3").Value = "Deprecate Vector3.  Replace it with position, direction and unit direction"
$ws.Range("C3").Value = 10

# Match formatting of the row below (row 5, the old row 3 "FSAA" row) for the two new rows.
$ws.Range("A2:C3").Font.Bold = $false

# Update selection to B2, matching the final saved workbook view.
$ws.Range("B2").Select()
